# Integrate freezeaccount / excelmanager output into the FrozenAccounts sheet:
# append the newly-frozen account as row 3 and make sure the sheet's
# right-to-left display flag is explicitly set to "off" (LTR), matching the
# worksheet's default reading order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Explicitly pin the sheet (and its window) to left-to-right display -
# mirrors the `rightToLeft="0"` attribute added to <sheetView>.
$ws.DisplayRightToLeft = $false
$excel.ActiveWindow.DisplayRightToLeft = $false

# Append the new frozen-account record as row 3.
$ws.Range("A3").Value = "3HkSLidfgeLyM1izEZMvB4eKHi94U4HWbBdfvY48Vpq3"
$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = "2fAbEYKeY2yfGAAhNCB2bpg1ACR3PvVACGha3Z6FM8HuosYzNUYVhzN7oPT2aqwWQsYMuTTzSwWaiK1YYgAcH3AZ"
